# Refactoring constrained alignment tree - more family-level branches
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 41-48: family column (D) changes from "Circoviridae" to "CRESS-1"
foreach ($r in 41..48) {
    $ws.Range("E$r").Copy()
    $ws.Range("D$r").PasteSpecial(-4122)
    $ws.Range("D$r").Value = "CRESS-1"
}

# Row 49: family column (D) changes from "Circoviridae" to "CRESS-2"
$ws.Range("E49").Copy()
$ws.Range("D49").PasteSpecial(-4122)
$ws.Range("D49").Value = "CRESS-2"

# Update the sheet view - scroll position and selection
$ws.Range("C47").Select()
$ws.Range("A1:H53").Select()
$excel.ActiveWindow.ScrollRow = 35
